# Add 2022-Q3 data:
#  - the existing "2021-Q1" detail sheet is renamed to "2022-Q3" and
#    re-populated with the new quarter's fund data
#  - a brand-new "2021-Q1" sheet is created (after "2022-Q3") holding an
#    exact copy of what used to be on the "2021-Q1" sheet
#  - the "总计" (totals) summary sheet gets a new row for 2022-Q3 and keeps
#    the old 2021-Q1 row (shifted down one row)

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$oldQ1   = $wb.Worksheets.Item(2)     # currently named "2021-Q1"

# ---------------------------------------------------------------------
# 1) Turn the original sheet (still holding 2021-Q1's data) into the new
#    2022-Q3 sheet by renaming the tab (its cells are overwritten below).
#    Renaming first frees up the "2021-Q1" name for the new sheet added
#    in step 2, and keeps this sheet's sheetId/r:id (it becomes rId2).
# ---------------------------------------------------------------------
$q3 = $oldQ1
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Create the new sheet that will keep the old "2021-Q1" detail data,
#    positioned right after 2022-Q3, and copy the original contents
#    (still sitting in $q3's cells) over verbatim before overwriting them.
# ---------------------------------------------------------------------
$newQ1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q3)
$newQ1.Name = "2021-Q1"

$q3.Range("B1:H2").Copy($newQ1.Range("B1"))
$q3.Range("A2").Copy($newQ1.Range("A2"))

# helper: write a value as TEXT (no auto number coercion), without
# leaving any stray formatting behind. Uses a scratch cell far outside
# the sheet's real data, formats it as Text, copies the value across
# with Paste Special (values only) and clears the scratch cell again.
function Set-TextCell($ws, $addr, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

# headers (row 1) - style copied from the summary sheet's header style
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"
$summary.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# row 2
$q3.Cells.Item(2, 1).Value = 0
$summary.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)
Set-TextCell $q3 "B2" "010695"
Set-TextCell $q3 "C2" "华夏磐益一年定期开放混合"
Set-TextCell $q3 "D2" "15.90"
Set-TextCell $q3 "E2" "99.95"
Set-TextCell $q3 "F2" "4.56"
Set-TextCell $q3 "G2" "0.7250"
$q3.Cells.Item(2, 8).Value = 3

# row 3
$q3.Cells.Item(3, 1).Value = 1
$summary.Range("A2").Copy()
$q3.Range("A3").PasteSpecial(-4122)
Set-TextCell $q3 "B3" "009837"
Set-TextCell $q3 "C3" "华夏磐锐一年定期开放混合A"
Set-TextCell $q3 "D3" "14.02"
Set-TextCell $q3 "E3" "94.15"
Set-TextCell $q3 "F3" "4.65"
Set-TextCell $q3 "G3" "0.6519"
$q3.Cells.Item(3, 8).Value = 3

# row 4
$q3.Cells.Item(4, 1).Value = 2
$summary.Range("A2").Copy()
$q3.Range("A4").PasteSpecial(-4122)
Set-TextCell $q3 "B4" "009838"
Set-TextCell $q3 "C4" "华夏磐锐一年定期开放混合C"
Set-TextCell $q3 "D4" "0.39"
Set-TextCell $q3 "E4" "94.15"
Set-TextCell $q3 "F4" "4.65"
Set-TextCell $q3 "G4" "0.0181"
$q3.Cells.Item(4, 8).Value = 3

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: row 2 now describes 2022-Q3, and a
#    new row 3 holds what used to be in row 2 (the 2021-Q1 totals).
# ---------------------------------------------------------------------
$summary.Cells.Item(3, 2).Value = "2021-Q1"
$summary.Cells.Item(3, 3).Value = 1
$summary.Cells.Item(3, 4).Value = 0.01
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$summary.Cells.Item(3, 1).Value = 1

$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 3
$summary.Cells.Item(2, 4).Value = 1.4
